$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 1029
$ws.Cells.Item(11, 9).Value = 1029
$ws.Cells.Item(11, 11).Value = 1029
$ws.Cells.Item(11, 13).Value = -889
$ws.Cells.Item(12, 8).Value = 3788311.8
$ws.Cells.Item(12, 9).Value = 5681930.5
$ws.Cells.Item(12, 11).Value = 5681930.5
$ws.Cells.Item(12, 13).Value = -5681760.5
$ws.Cells.Item(28, 8).Value = 1025.909
$ws.Cells.Item(28, 9).Value = 339.84
$ws.Cells.Item(28, 11).Value = 339.84
$ws.Cells.Item(28, 13).Value = 145.16
$ws.Cells.Item(86, 8).Value = 6335.5264
$ws.Cells.Item(86, 9).Value = 5374.6665
$ws.Cells.Item(86, 10).Value = 7200.3
$ws.Cells.Item(86, 11).Value = 5374.6665
$ws.Cells.Item(86, 12).Value = 7200.3
$ws.Cells.Item(86, 13).Value = -4251.6665
$ws.Cells.Item(86, 14).Value = -9446.299999999999
$ws.Cells.Item(89, 8).Value = 6335.5264
$ws.Cells.Item(89, 9).Value = 5374.6665
$ws.Cells.Item(89, 10).Value = 7200.3
$ws.Cells.Item(89, 11).Value = 26873.3325
$ws.Cells.Item(89, 12).Value = 36001.5
$ws.Cells.Item(89, 13).Value = -21257.3325
$ws.Cells.Item(89, 14).Value = -47233.5
$ws.Cells.Item(113, 8).Value = 4502.2812
$ws.Cells.Item(113, 10).Value = 4202.6113
$ws.Cells.Item(113, 12).Value = 4202.6113
$ws.Cells.Item(113, 14).Value = -10710.6113
$ws.Cells.Item(137, 8).Value = 46651.8
$ws.Cells.Item(137, 9).Value = 82491.63
$ws.Cells.Item(137, 10).Value = 2847.5557
$ws.Cells.Item(137, 11).Value = 247474.89
$ws.Cells.Item(137, 12).Value = 8542.667099999999
$ws.Cells.Item(137, 13).Value = -244924.89
$ws.Cells.Item(137, 14).Value = -13642.6671
$ws.Cells.Item(138, 8).Value = 3658.2896
$ws.Cells.Item(138, 9).Value = 2934.8
$ws.Cells.Item(138, 11).Value = 8804.400000000001
$ws.Cells.Item(138, 13).Value = -3664.400000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1048193.56
$ws.Cells.Item(2, 9).Value = 1131862.8
$ws.Cells.Item(2, 10).Value = 2328.25
$ws.Cells.Item(2, 11).Value = 1131862.8
$ws.Cells.Item(2, 12).Value = 2328.25
$ws.Cells.Item(2, 13).Value = -1131749.8
$ws.Cells.Item(2, 14).Value = -2554.25
$ws.Cells.Item(45, 8).Value = 5293773.5
$ws.Cells.Item(45, 9).Value = 6804176
$ws.Cells.Item(45, 11).Value = 6804176
$ws.Cells.Item(45, 13).Value = -6803799
$ws.Cells.Item(61, 8).Value = 4381.0527
$ws.Cells.Item(61, 9).Value = 4885.909
$ws.Cells.Item(61, 10).Value = 3686.875
$ws.Cells.Item(61, 11).Value = 4885.909
$ws.Cells.Item(61, 12).Value = 3686.875
$ws.Cells.Item(61, 13).Value = -4673.909
$ws.Cells.Item(61, 14).Value = -4110.875
$ws.Cells.Item(74, 8).Value = 33895.38
$ws.Cells.Item(74, 9).Value = 2370.5652
$ws.Cells.Item(74, 11).Value = 2370.5652
$ws.Cells.Item(74, 13).Value = -1496.5652
$ws.Cells.Item(77, 8).Value = 33895.38
$ws.Cells.Item(77, 9).Value = 2370.5652
$ws.Cells.Item(77, 11).Value = 11852.826
$ws.Cells.Item(77, 13).Value = -7484.826000000001
$ws.Cells.Item(110, 8).Value = 1158830.8
$ws.Cells.Item(110, 9).Value = 1853002.8
$ws.Cells.Item(110, 10).Value = 1877.4445
$ws.Cells.Item(110, 11).Value = 1853002.8
$ws.Cells.Item(110, 12).Value = 1877.4445
$ws.Cells.Item(110, 13).Value = -1850957.8
$ws.Cells.Item(110, 14).Value = -5967.4445
$ws.Cells.Item(114, 8).Value = 82549.664
$ws.Cells.Item(114, 10).Value = 82549.664
$ws.Cells.Item(114, 12).Value = 82549.664
$ws.Cells.Item(114, 14).Value = -91227.664
$ws.Cells.Item(116, 8).Value = 1048193.56
$ws.Cells.Item(116, 9).Value = 1131862.8
$ws.Cells.Item(116, 10).Value = 2328.25
$ws.Cells.Item(116, 11).Value = 1131862.8
$ws.Cells.Item(116, 12).Value = 2328.25
$ws.Cells.Item(116, 13).Value = -1129568.8
$ws.Cells.Item(116, 14).Value = -6916.25
$ws.Cells.Item(122, 8).Value = 1635117
$ws.Cells.Item(122, 9).Value = 1881049
$ws.Cells.Item(122, 10).Value = 1230052.5
$ws.Cells.Item(122, 11).Value = 5643147
$ws.Cells.Item(122, 12).Value = 3690157.5
$ws.Cells.Item(122, 13).Value = -5640697
$ws.Cells.Item(122, 14).Value = -3695057.5
$ws.Cells.Item(132, 8).Value = 3860.4546
$ws.Cells.Item(132, 10).Value = 6748.5
$ws.Cells.Item(132, 12).Value = 20245.5
$ws.Cells.Item(132, 14).Value = -25305.5
$ws.Cells.Item(136, 8).Value = 4381.0527
$ws.Cells.Item(136, 9).Value = 4885.909
$ws.Cells.Item(136, 10).Value = 3686.875
$ws.Cells.Item(136, 11).Value = 14657.727
$ws.Cells.Item(136, 12).Value = 11060.625
$ws.Cells.Item(136, 13).Value = -12107.727
$ws.Cells.Item(136, 14).Value = -16160.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1048193.56
$ws.Cells.Item(3, 9).Value = 1131862.8
$ws.Cells.Item(3, 10).Value = 2328.25
$ws.Cells.Item(3, 11).Value = 1131862.8
$ws.Cells.Item(3, 12).Value = 2328.25
$ws.Cells.Item(3, 13).Value = -1131748.8
$ws.Cells.Item(3, 14).Value = -2556.25
$ws.Cells.Item(86, 8).Value = 2860130.8
$ws.Cells.Item(86, 9).Value = 4169301.5
$ws.Cells.Item(86, 11).Value = 4169301.5
$ws.Cells.Item(86, 13).Value = -4168178.5
$ws.Cells.Item(89, 8).Value = 2860130.8
$ws.Cells.Item(89, 9).Value = 4169301.5
$ws.Cells.Item(89, 11).Value = 20846507.5
$ws.Cells.Item(89, 13).Value = -20840891.5
$ws.Cells.Item(105, 8).Value = 3908269.2
$ws.Cells.Item(105, 9).Value = 4809808.5
$ws.Cells.Item(105, 11).Value = 4809808.5
$ws.Cells.Item(105, 13).Value = -4808061.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 436.63635
$ws.Cells.Item(22, 9).Value = 422.66666
$ws.Cells.Item(22, 11).Value = 422.66666
$ws.Cells.Item(22, 13).Value = -72.66665999999998
$ws.Cells.Item(58, 8).Value = 6085.5674
$ws.Cells.Item(58, 9).Value = 7820.174
$ws.Cells.Item(58, 11).Value = 7820.174
$ws.Cells.Item(58, 13).Value = -7617.174
$ws.Cells.Item(107, 8).Value = 1339.0209
$ws.Cells.Item(107, 9).Value = 1379.6757
$ws.Cells.Item(107, 11).Value = 1379.6757
$ws.Cells.Item(107, 13).Value = 540.3243
$ws.Cells.Item(134, 8).Value = 1852.7273
$ws.Cells.Item(134, 9).Value = 1211.8966
$ws.Cells.Item(134, 11).Value = 3635.6898
$ws.Cells.Item(134, 13).Value = -1100.6898
$ws.Cells.Item(136, 8).Value = 6085.5674
$ws.Cells.Item(136, 9).Value = 7820.174
$ws.Cells.Item(136, 11).Value = 23460.522
$ws.Cells.Item(136, 13).Value = -20910.522

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 36277320
$ws.Cells.Item(4, 10).Value = 10319.583
$ws.Cells.Item(4, 12).Value = 30958.749
$ws.Cells.Item(4, 14).Value = -31182.749
$ws.Cells.Item(47, 8).Value = 965
$ws.Cells.Item(47, 9).Value = 151
$ws.Cells.Item(47, 11).Value = 453
$ws.Cells.Item(47, 13).Value = -22
$ws.Cells.Item(107, 8).Value = 2232.2856
$ws.Cells.Item(107, 9).Value = 2744.25
$ws.Cells.Item(107, 11).Value = 8232.75
$ws.Cells.Item(107, 13).Value = -6312.75
$ws.Cells.Item(113, 8).Value = 3705.3809
$ws.Cells.Item(113, 10).Value = 1896.8182
$ws.Cells.Item(113, 12).Value = 5690.4546
$ws.Cells.Item(113, 14).Value = -10030.4546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 29154838
$ws.Cells.Item(80, 10).Value = 26513.75
$ws.Cells.Item(80, 12).Value = 26513.75
$ws.Cells.Item(80, 14).Value = -28509.75
$ws.Cells.Item(83, 8).Value = 29154838
$ws.Cells.Item(83, 10).Value = 26513.75
$ws.Cells.Item(83, 12).Value = 132568.75
$ws.Cells.Item(83, 14).Value = -142552.75
$ws.Cells.Item(113, 8).Value = 4168734
$ws.Cells.Item(113, 9).Value = 5209834.5
$ws.Cells.Item(113, 11).Value = 5209834.5
$ws.Cells.Item(113, 13).Value = -5207664.5
$ws.Cells.Item(132, 8).Value = 2929.4866
$ws.Cells.Item(132, 9).Value = 2809.25
$ws.Cells.Item(132, 10).Value = 3303.5557
$ws.Cells.Item(132, 11).Value = 8427.75
$ws.Cells.Item(132, 12).Value = 9910.667099999999
$ws.Cells.Item(132, 13).Value = -5897.75
$ws.Cells.Item(132, 14).Value = -14970.6671
$ws.Cells.Item(141, 8).Value = 54183
$ws.Cells.Item(141, 10).Value = 65856.2
$ws.Cells.Item(141, 12).Value = 65856.2
$ws.Cells.Item(141, 14).Value = -76216.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 13898779
$ws.Cells.Item(93, 9).Value = 17545168
$ws.Cells.Item(93, 11).Value = 17545168
$ws.Cells.Item(93, 13).Value = -17543920
$ws.Cells.Item(132, 8).Value = 13687.5625
$ws.Cells.Item(132, 9).Value = 14538.538
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 43615.614
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).Value = -41085.614
$ws.Cells.Item(132, 14).Value = -35060
$ws.Cells.Item(136, 8).Value = 51056.684
$ws.Cells.Item(136, 9).Value = 158170.53
$ws.Cells.Item(136, 10).Value = 6137.968
$ws.Cells.Item(136, 11).Value = 474511.59
$ws.Cells.Item(136, 12).Value = 18413.904
$ws.Cells.Item(136, 13).Value = -471961.59
$ws.Cells.Item(136, 14).Value = -23513.904

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1190.3
$ws.Cells.Item(100, 9).Value = 1337.2916
$ws.Cells.Item(100, 11).Value = 2674.5832
$ws.Cells.Item(100, 13).Value = -2133.5832
$ws.Cells.Item(107, 8).Value = 32259366
$ws.Cells.Item(107, 9).Value = 58824160
$ws.Cells.Item(107, 11).Value = 176472480
$ws.Cells.Item(107, 13).Value = -176470560
$ws.Cells.Item(126, 8).Value = 4080.7856
$ws.Cells.Item(126, 10).Value = 6106.6665
$ws.Cells.Item(126, 12).Value = 18319.9995
$ws.Cells.Item(126, 14).Value = -23259.9995
$ws.Cells.Item(132, 8).Value = 22979228
$ws.Cells.Item(132, 9).Value = 31255400
$ws.Cells.Item(132, 10).Value = 909434.8
$ws.Cells.Item(132, 11).Value = 93766200
$ws.Cells.Item(132, 12).Value = 2728304.4
$ws.Cells.Item(132, 13).Value = -93763670
$ws.Cells.Item(132, 14).Value = -2733364.4
$ws.Cells.Item(136, 8).Value = 884.53125
$ws.Cells.Item(136, 9).Value = 795.3509
$ws.Cells.Item(136, 10).Value = 1610.7142
$ws.Cells.Item(136, 11).Value = 2386.0527
$ws.Cells.Item(136, 12).Value = 4832.142599999999
$ws.Cells.Item(136, 14).Value = -9932.142599999999
